$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") bumps from 46063 to 46064 for rows 2-14
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = 46064
}

# Rows 7-14: data rows were re-sorted/refreshed, changing Beteckning (A),
# Datum (B) and Area (G) values per the new source order.
$rowData = @(
    @{ Row = 7;  A = "A 14271-2021"; B = 44278;              G = 6.7 },
    @{ Row = 8;  A = "A 62884-2021"; B = 44504;              G = 0.8 },
    @{ Row = 9;  A = "A 25634-2025"; B = 45803.59570601852;  G = 6 },
    @{ Row = 10; A = "A 28266-2025"; B = 45818.56381944445;  G = 1.9 },
    @{ Row = 11; A = "A 25015-2023"; B = 45085.6989699074;   G = 1.8 },
    @{ Row = 12; A = "A 19922-2025"; B = 45771.63034722222;  G = 10.1 },
    @{ Row = 13; A = "A 60024-2025"; B = 45992;              G = 1.1 },
    @{ Row = 14; A = "A 3402-2026";  B = 46042.39047453704;  G = 5.5 }
)

foreach ($item in $rowData) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 7).Value = $item.G
}
